$wb = $excel.ActiveWorkbook

$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: $newVersion"
$wsAbout.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Tashan Coal Mine, China, M0358, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 15; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
